# Auto-generated edit script: refreshes cached market-price-derived
# columns (currentAveragePrice / currentAveragePriceNQ / currentAveragePriceHQ /
# LevePriceNQ / LevePriceHQ / LeveProfitNQ / LeveProfitHQ) on a handful of rows
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets, as produced by the scheduled
# market-data refresh runner.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC row 20
$ws.Range("H20").Value = 957.75
$ws.Range("I20").Value = 957.75
$ws.Range("K20").Value = 957.75
$ws.Range("M20").Value = -727.75
# ALC row 35
$ws.Range("H35").Value = 957.75
$ws.Range("I35").Value = 957.75
$ws.Range("K35").Value = 957.75
$ws.Range("M35").Value = -578.75
# ALC row 40
$ws.Range("H40").Value = 3275.2424
$ws.Range("I40").Value = 2572.1904
$ws.Range("J40").Value = 4505.5835
$ws.Range("K40").Value = 2572.1904
$ws.Range("L40").Value = 4505.5835
$ws.Range("M40").Value = -2397.1904
$ws.Range("N40").Value = -4855.5835
# ALC row 86
$ws.Range("H86").Value = 4390397.5
$ws.Range("I86").Value = 4199.6665
$ws.Range("J86").Value = 8776596
$ws.Range("K86").Value = 4199.6665
$ws.Range("L86").Value = 8776596
$ws.Range("M86").Value = -3076.6665
$ws.Range("N86").Value = -8778842
# ALC row 89
$ws.Range("H89").Value = 4390397.5
$ws.Range("I89").Value = 4199.6665
$ws.Range("J89").Value = 8776596
$ws.Range("K89").Value = 20998.3325
$ws.Range("L89").Value = 43882980
$ws.Range("M89").Value = -15382.3325
$ws.Range("N89").Value = -43894212
# ALC row 113
$ws.Range("H113").Value = 8836.083000000001
$ws.Range("I113").Value = 6900
$ws.Range("K113").Value = 6900
$ws.Range("M113").Value = -3646
# ALC row 116
$ws.Range("H116").Value = 4999.5
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = 4999.5
$ws.Range("K116").Value = 0
$ws.Range("L116").Value = 4999.5
$ws.Range("M116").ClearContents()
$ws.Range("N116").Value = -11883.5
# ALC row 132
$ws.Range("H132").Value = 37441.5
$ws.Range("I132").Value = 1786.8889
$ws.Range("J132").Value = 101619.8
$ws.Range("K132").Value = 5360.6667
$ws.Range("L132").Value = 304859.4
$ws.Range("M132").Value = -2830.6667
$ws.Range("N132").Value = -309919.4
# ALC row 137
$ws.Range("H137").Value = 3870.7742
$ws.Range("I137").Value = 7828
$ws.Range("J137").Value = 1986.381
$ws.Range("K137").Value = 23484
$ws.Range("L137").Value = 5959.143
$ws.Range("M137").Value = -20934
$ws.Range("N137").Value = -11059.143

$ws = $wb.Worksheets.Item("ARM")
# ARM row 31
$ws.Range("H31").Value = 21873.8
$ws.Range("I31").Value = 18842.25
$ws.Range("J31").Value = 34000
$ws.Range("K31").Value = 18842.25
$ws.Range("L31").Value = 34000
$ws.Range("M31").Value = -18548.25
$ws.Range("N31").Value = -34588
# ARM row 32
$ws.Range("H32").Value = 23817800
$ws.Range("I32").Value = 23817800
$ws.Range("K32").Value = 23817800
$ws.Range("M32").Value = -23817513
# ARM row 33
$ws.Range("H33").Value = 3693
$ws.Range("I33").Value = 3693
$ws.Range("K33").Value = 3693
$ws.Range("M33").Value = -3364
# ARM row 61
$ws.Range("H61").Value = 5407.4287
$ws.Range("I61").Value = 4125.625
$ws.Range("J61").Value = 7116.5
$ws.Range("K61").Value = 4125.625
$ws.Range("L61").Value = 7116.5
$ws.Range("M61").Value = -3913.625
$ws.Range("N61").Value = -7540.5
# ARM row 74
$ws.Range("H74").Value = 3138.5454
$ws.Range("I74").Value = 3367.7144
$ws.Range("K74").Value = 3367.7144
$ws.Range("M74").Value = -2493.7144
# ARM row 77
$ws.Range("H77").Value = 3138.5454
$ws.Range("I77").Value = 3367.7144
$ws.Range("K77").Value = 16838.572
$ws.Range("M77").Value = -12470.572
# ARM row 102
$ws.Range("H102").Value = 66669810
$ws.Range("I102").Value = 4570
$ws.Range("K102").Value = 4570
$ws.Range("M102").Value = -2948
# ARM row 110
$ws.Range("H110").Value = 849.9
$ws.Range("I110").Value = 844.3333
$ws.Range("J110").Value = 900
$ws.Range("K110").Value = 844.3333
$ws.Range("L110").Value = 900
$ws.Range("M110").Value = 1200.6667
$ws.Range("N110").Value = -4990
# ARM row 136
$ws.Range("H136").Value = 5407.4287
$ws.Range("I136").Value = 4125.625
$ws.Range("J136").Value = 7116.5
$ws.Range("K136").Value = 12376.875
$ws.Range("L136").Value = 21349.5
$ws.Range("M136").Value = -9826.875
$ws.Range("N136").Value = -26449.5

$ws = $wb.Worksheets.Item("BSM")
# BSM row 20
$ws.Range("H20").Value = 5455.7144
$ws.Range("I20").Value = 5839
$ws.Range("J20").Value = 4497.5
$ws.Range("K20").Value = 5839
$ws.Range("L20").Value = 4497.5
$ws.Range("M20").Value = -5592
$ws.Range("N20").Value = -4991.5
# BSM row 22
$ws.Range("H22").Value = 3137.5
$ws.Range("I22").Value = 3137.5
$ws.Range("K22").Value = 3137.5
$ws.Range("M22").Value = -2964.5
# BSM row 31
$ws.Range("H31").Value = 20000
$ws.Range("J31").Value = 20000
$ws.Range("L31").Value = 20000
$ws.Range("N31").Value = -20504
# BSM row 105
$ws.Range("H105").Value = 4008.3157
$ws.Range("I105").Value = 3238.25
$ws.Range("K105").Value = 3238.25
$ws.Range("M105").Value = -1491.25
# BSM row 134
$ws.Range("H134").Value = 2705.3137
$ws.Range("I134").Value = 1748.6097
$ws.Range("K134").Value = 5245.8291
$ws.Range("M134").Value = -2710.8291

$ws = $wb.Worksheets.Item("CRP")
# CRP row 2
$ws.Range("H2").Value = 122
$ws.Range("I2").Value = 149.33333
$ws.Range("J2").Value = 40
$ws.Range("K2").Value = 149.33333
$ws.Range("L2").Value = 40
$ws.Range("M2").Value = -36.33332999999999
$ws.Range("N2").Value = -266
# CRP row 58
$ws.Range("H58").Value = 5958.154
$ws.Range("I58").Value = 4363.5
$ws.Range("K58").Value = 4363.5
$ws.Range("M58").Value = -4160.5
# CRP row 105
$ws.Range("H105").Value = 50017900
$ws.Range("I105").Value = 0
$ws.Range("K105").Value = 0
$ws.Range("M105").ClearContents()
# CRP row 136
$ws.Range("H136").Value = 5958.154
$ws.Range("I136").Value = 4363.5
$ws.Range("K136").Value = 13090.5
$ws.Range("M136").Value = -10540.5

$ws = $wb.Worksheets.Item("CUL")
# CUL row 107
$ws.Range("H107").Value = 248.81818
$ws.Range("I107").Value = 148.75
$ws.Range("J107").Value = 306
$ws.Range("K107").Value = 446.25
$ws.Range("L107").Value = 918
$ws.Range("M107").Value = 1473.75
$ws.Range("N107").Value = -4758
# CUL row 113
$ws.Range("H113").Value = 871.13043
$ws.Range("J113").Value = 940.8
$ws.Range("L113").Value = 2822.4
$ws.Range("N113").Value = -7162.4
# CUL row 122
$ws.Range("H122").Value = 1043.9412
$ws.Range("I122").Value = 962.875
$ws.Range("J122").Value = 1116
$ws.Range("K122").Value = 8665.875
$ws.Range("L122").Value = 10044
$ws.Range("M122").Value = -6215.875
$ws.Range("N122").Value = -14944

$ws = $wb.Worksheets.Item("GSM")
# GSM row 80
$ws.Range("H80").Value = 22299968
$ws.Range("I80").Value = 189050.67
$ws.Range("K80").Value = 189050.67
$ws.Range("M80").Value = -188052.67
# GSM row 83
$ws.Range("H83").Value = 22299968
$ws.Range("I83").Value = 189050.67
$ws.Range("K83").Value = 945253.3500000001
$ws.Range("M83").Value = -940261.3500000001
# GSM row 126
$ws.Range("H126").Value = 4622
$ws.Range("I126").Value = 1933
$ws.Range("K126").Value = 5799
$ws.Range("M126").Value = -3329
# GSM row 132
$ws.Range("H132").Value = 3024.8572
$ws.Range("I132").Value = 2554.8
$ws.Range("K132").Value = 7664.400000000001
$ws.Range("M132").Value = -5134.400000000001
# GSM row 135
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()
# GSM row 140
$ws.Range("H140").Value = 80000
$ws.Range("J140").Value = 80000
$ws.Range("L140").Value = 80000
$ws.Range("N140").Value = -90360

$ws = $wb.Worksheets.Item("LTW")
# LTW row 7
$ws.Range("H7").Value = 4676
$ws.Range("I7").Value = 412
$ws.Range("J7").Value = 5332
$ws.Range("K7").Value = 412
$ws.Range("L7").Value = 5332
$ws.Range("M7").Value = -300
$ws.Range("N7").Value = -5556
# LTW row 68
$ws.Range("H68").Value = 211419.83
$ws.Range("I68").Value = 127370.81
$ws.Range("K68").Value = 127370.81
$ws.Range("M68").Value = -126621.81
# LTW row 71
$ws.Range("H71").Value = 211419.83
$ws.Range("I71").Value = 127370.81
$ws.Range("K71").Value = 636854.05
$ws.Range("M71").Value = -633110.05
# LTW row 122
$ws.Range("H122").Value = 5239.5
$ws.Range("I122").Value = 2685
$ws.Range("J122").Value = 6942.5
$ws.Range("K122").Value = 8055
$ws.Range("L122").Value = 20827.5
$ws.Range("M122").Value = -5605
$ws.Range("N122").Value = -25727.5
# LTW row 126
$ws.Range("H126").Value = 4676
$ws.Range("I126").Value = 412
$ws.Range("J126").Value = 5332
$ws.Range("K126").Value = 1236
$ws.Range("L126").Value = 15996
$ws.Range("M126").Value = 1234
$ws.Range("N126").Value = -20936
# LTW row 132
$ws.Range("H132").Value = 4840.4185
$ws.Range("I132").Value = 3386.2646
$ws.Range("J132").Value = 10333.889
$ws.Range("K132").Value = 10158.7938
$ws.Range("L132").Value = 31001.667
$ws.Range("M132").Value = -7628.793799999999
$ws.Range("N132").Value = -36061.667
# LTW row 136
$ws.Range("H136").Value = 4832.273
$ws.Range("I136").Value = 2877.7334
$ws.Range("J136").Value = 6461.0557
$ws.Range("K136").Value = 8633.200199999999
$ws.Range("L136").Value = 19383.1671
$ws.Range("M136").Value = -6083.200199999999
$ws.Range("N136").Value = -24483.1671

$ws = $wb.Worksheets.Item("WVR")
# WVR row 116
$ws.Range("H116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("N116").ClearContents()
# WVR row 132
$ws.Range("H132").Value = 1889.0834
$ws.Range("I132").Value = 1013.53845
$ws.Range("J132").Value = 4165.5
$ws.Range("K132").Value = 3040.61535
$ws.Range("L132").Value = 12496.5
$ws.Range("M132").Value = -510.61535
$ws.Range("N132").Value = -17556.5

Write-Output "Updated market-price columns on $( 45 ) rows across $( 8 ) sheets."